$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) contains numeric-looking text (e.g. trailing
# zeros like "1.00", or thousand-dot separated values like "61.112.33")
# that must stay as literal text rather than being auto-converted to
# real numbers by Excel, so force the column to Text format first.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "61.112.33"
$ws.Range("E2").Value = "  -0.65%  "

# Row 3
$ws.Range("D3").Value = "2.949.61"
$ws.Range("E3").Value = "  -1.44%  "

# Row 4
$ws.Range("E4").Value = "  +0.04%  "

# Row 5
$ws.Range("D5").Value = "535.55"
$ws.Range("E5").Value = "  -0.24%  "

# Row 6
$ws.Range("D6").Value = "132.84"
$ws.Range("E6").Value = "  +0.11%  "

# Row 7
$ws.Range("E7").Value = "  +0.01%  "

# Row 8
$ws.Range("D8").Value = "2.948.44"
$ws.Range("E8").Value = "  -1.40%  "

# Row 9
$ws.Range("D9").Value = "0.481"
$ws.Range("E9").Value = "  -2.51%  "

# Row 10
$ws.Range("D10").Value = "6.33"
$ws.Range("E10").Value = "  +3.97%  "

# Row 11
$ws.Range("D11").Value = "0.144"
$ws.Range("E11").Value = "  -2.15%  "

# Row 12
$ws.Range("D12").Value = "0.436"
$ws.Range("E12").Value = "  -2.25%  "

# Row 13
$ws.Range("D13").Value = "0.0000216"
$ws.Range("E13").Value = "  -2.50%  "

# Row 14
$ws.Range("D14").Value = "32.96"
$ws.Range("E14").Value = "  -2.33%  "

# Row 15
$ws.Range("D15").Value = "3.459.42"
$ws.Range("E15").Value = "  -0.50%  "

# Row 16
$ws.Range("D16").Value = "61.169.12"
$ws.Range("E16").Value = "  -0.57%  "

# Row 17
$ws.Range("B17").Value = "TRON"
$ws.Range("C17").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D17").Value = "0.107"
$ws.Range("E17").Value = "  -2.13%  "

# Row 18
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "2.952.49"
$ws.Range("E18").Value = "  -1.36%  "

# Row 19
$ws.Range("D19").Value = "6.47"
$ws.Range("E19").Value = "  -1.93%  "

# Row 20
$ws.Range("D20").Value = "460.37"
$ws.Range("E20").Value = "  -1.50%  "

# Row 21
$ws.Range("D21").Value = "13.16"
$ws.Range("E21").Value = "  +0.28%  "

# Row 22
$ws.Range("D22").Value = "0.639"
$ws.Range("E22").Value = "  -4.34%  "

# Row 23
$ws.Range("D23").Value = "6.90"
$ws.Range("E23").Value = "  -0.26%  "

# Row 24
$ws.Range("D24").Value = "78.16"
$ws.Range("E24").Value = "  -2.33%  "

# Row 25
$ws.Range("D25").Value = "12.19"
$ws.Range("E25").Value = "  +2.15%  "

# Row 26
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  +0.12%  "

# Row 27
$ws.Range("D27").Value = "2.66"
$ws.Range("E27").Value = "  -0.78%  "

# Row 28
$ws.Range("B28").Value = "FirstDigitalUSD"
$ws.Range("C28").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D28").Value = "0.999"
$ws.Range("E28").Value = "  -0.17%  "

# Row 29
$ws.Range("B29").Value = "RenderToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D29").Value = "7.40"
$ws.Range("E29").Value = "  -4.33%  "

# Row 30
$ws.Range("B30").Value = "ImmutableX"
$ws.Range("C30").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D30").Value = "1.94"
$ws.Range("E30").Value = "  +3.22%  "

# Row 31
$ws.Range("D31").Value = "25.02"
$ws.Range("E31").Value = "  -1.81%  "

# Row 32
$ws.Range("D32").Value = "1.11"
$ws.Range("E32").Value = "  -2.82%  "

# Row 33
$ws.Range("D33").Value = "2.30"
$ws.Range("E33").Value = "  +0.34%  "

# Row 34
$ws.Range("D34").Value = "5.45"
$ws.Range("E34").Value = "  +0.17%  "

# Row 35
$ws.Range("D35").Value = "53.86"
$ws.Range("E35").Value = "  -3.02%  "

# Row 36
$ws.Range("D36").Value = "5.72"
$ws.Range("E36").Value = "  -2.46%  "

# Row 37
$ws.Range("D37").Value = "443.42"
$ws.Range("E37").Value = "  -3.62%  "

# Row 38
$ws.Range("D38").Value = "0.0784"
$ws.Range("E38").Value = "  -0.56%  "

# Row 39
$ws.Range("D39").Value = "0.0382"
$ws.Range("E39").Value = "  +0.34%  "

# Row 40
$ws.Range("D40").Value = "2.899.26"
$ws.Range("E40").Value = "  -8.40%  "

# Row 41
$ws.Range("D41").Value = "0.114"
$ws.Range("E41").Value = "  -4.98%  "

# Row 42
$ws.Range("D42").Value = "7.84"
$ws.Range("E42").Value = "  -2.80%  "

# Row 43
$ws.Range("B43").Value = "USDe"
$ws.Range("C43").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D43").Value = "1.00"
$ws.Range("E43").Value = "  -0.07%  "

# Row 44
$ws.Range("D44").Value = "26.03"
$ws.Range("E44").Value = "  +1.22%  "

# Row 45
$ws.Range("B45").Value = "dogwifhat"
$ws.Range("C45").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D45").Value = "2.34"
$ws.Range("E45").Value = "  -2.52%  "

# Row 46
$ws.Range("B46").Value = "TheGraph"
$ws.Range("C46").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D46").Value = "0.242"
$ws.Range("E46").Value = "  -0.03%  "

# Row 47
$ws.Range("B47").Value = "Stellar"
$ws.Range("C47").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D47").Value = "0.107"
$ws.Range("E47").Value = "  -0.37%  "

# Row 48
$ws.Range("B48").Value = "Fetch.AI"
$ws.Range("C48").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D48").Value = "1.93"
$ws.Range("E48").Value = "  -2.94%  "

# Row 49
$ws.Range("D49").Value = "113.55"
$ws.Range("E49").Value = "  -3.41%  "

# Row 50
$ws.Range("D50").Value = "0.0₃0479"
$ws.Range("E50").Value = "  -2.20%  "

# Row 51
$ws.Range("D51").Value = "1.21"
$ws.Range("E51").Value = "  -2.69%  "
